$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update "Förändrad" (changed) date column C for rows 2..347 from 45189 to 45190
$ws.Range("C2:C347").Value = 45190

# 2. Set row 347 explicit height (ht="15" customHeight="1")
$ws.Rows.Item(347).RowHeight = 15

# 3. Copy formatting (styles) from row 347 down into the two new rows 348 and 349
#    so that date columns (B,C) keep the date style and column R keeps its style.
$ws.Range("A347:R347").Copy()
$ws.Range("A348:R348").PasteSpecial(-4122)
$ws.Range("A349:R349").PasteSpecial(-4122)
$ws.Rows.Item(348).RowHeight = 15
# Column F (Markägare) is always blank/unset in this sheet - clear it so it
# doesn't linger as an empty styled cell.
$ws.Cells.Item(348, 6).Clear()
$ws.Cells.Item(349, 6).Clear()

# 4. Populate row 348
$ws.Cells.Item(348, 1).Value = "A 44586-2023"
$ws.Cells.Item(348, 2).Value = 45189
$ws.Cells.Item(348, 3).Value = 45190
$ws.Cells.Item(348, 4).Value = "VÄSTMANLANDS LÄN"
$ws.Cells.Item(348, 5).Value = "KÖPING"
$ws.Cells.Item(348, 7).Value = 3.2
$ws.Cells.Item(348, 8).Value = 0
$ws.Cells.Item(348, 9).Value = 0
$ws.Cells.Item(348, 10).Value = 0
$ws.Cells.Item(348, 11).Value = 0
$ws.Cells.Item(348, 12).Value = 0
$ws.Cells.Item(348, 13).Value = 0
$ws.Cells.Item(348, 14).Value = 0
$ws.Cells.Item(348, 15).Value = 0
$ws.Cells.Item(348, 16).Value = 0
$ws.Cells.Item(348, 17).Value = 0

# 5. Populate row 349
$ws.Cells.Item(349, 1).Value = "A 44582-2023"
$ws.Cells.Item(349, 2).Value = 45189
$ws.Cells.Item(349, 3).Value = 45190
$ws.Cells.Item(349, 4).Value = "VÄSTMANLANDS LÄN"
$ws.Cells.Item(349, 5).Value = "KÖPING"
$ws.Cells.Item(349, 7).Value = 1.7
$ws.Cells.Item(349, 8).Value = 0
$ws.Cells.Item(349, 9).Value = 0
$ws.Cells.Item(349, 10).Value = 0
$ws.Cells.Item(349, 11).Value = 0
$ws.Cells.Item(349, 12).Value = 0
$ws.Cells.Item(349, 13).Value = 0
$ws.Cells.Item(349, 14).Value = 0
$ws.Cells.Item(349, 15).Value = 0
$ws.Cells.Item(349, 16).Value = 0
$ws.Cells.Item(349, 17).Value = 0

$excel.CutCopyMode = $false
